# Apply updated crypto price/volume figures to the worksheet (matches the
# upstream "Updated cryptos list" GitHub Actions commit).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell -> new text value. Values that look like plain decimals (e.g. "217.82")
# are prefixed with a leading single-quote so Excel stores them as literal text
# (matching the original inlineStr cells) instead of auto-converting them to numbers.
$updates = [ordered]@{
    'D2' = '26.743.92'
    'E2' = '  +0.50%  '
    'D3' = '1.640.90'
    'E3' = '  -0.01%  '
    'E4' = '  +0.25%  '
    'D5' = '''217.82'
    'E5' = '  +1.41%  '
    'D6' = '''0.504'
    'E6' = '  -0.05%  '
    'E7' = '  +0.41%  '
    'E8' = '  +0.28%  '
    'E9' = '  +0.03%  '
    'D10' = '''19.11'
    'E10' = '  +0.11%  '
    'E11' = '  +0.21%  '
    'D12' = '1.869.78'
    'E12' = '  -0.01%  '
    'D13' = '1.642.53'
    'E13' = '  +0.20%  '
    'E14' = '  -0.37%  '
    'E15' = '  -0.33%  '
    'D16' = '''64.68'
    'E16' = '  -0.22%  '
    'D17' = '26.734.32'
    'E18' = '  -1.07%  '
    'D19' = '''214.16'
    'E19' = '  -0.56%  '
    'E20' = '  +0.32%  '
    'E21' = '  +0.93%  '
    'D22' = '''2.37'
    'E22' = '  +7.92%  '
    'E23' = '  -0.44%  '
    'E24' = '  -1.97%  '
    'D25' = '''145.63'
    'E25' = '  +0.40%  '
    'E26' = '  +0.25%  '
    'E27' = '  -0.80%  '
    'E28' = '  +0.67%  '
    'D29' = '''15.66'
    'E29' = '  -0.13%  '
    'E30' = '  -0.66%  '
    'E31' = '  +1.37%  '
    'E32' = '  +0.98%  '
    'D33' = '''3.02'
    'E33' = '  +0.78%  '
    'D34' = '1.287.36'
    'E34' = '  +0.83%  '
    'D35' = '''1.54'
    'E35' = '  -0.12%  '
    'E36' = '  +1.29%  '
    'E37' = '  -0.19%  '
    'D38' = '''0.537'
    'E38' = '  +1.13%  '
    'D39' = '''0.816'
    'E39' = '  -0.80%  '
    'E40' = '  +0.51%  '
    'E41' = '  -0.53%  '
    'E42' = '  -1.30%  '
    'E43' = '  -2.36%  '
    'D44' = '1.779.48'
    'E44' = '  -0.01%  '
    'D45' = '''61.02'
    'E45' = '  +3.20%  '
    'D46' = '''91.68'
    'E46' = '  +0.17%  '
    'E47' = '  -0.01%  '
    'E48' = '  +0.40%  '
    'D49' = '''7.62'
    'E49' = '  -1.67%  '
    'E50' = '  +0.42%  '
    'D51' = '''0.406'
    'E51' = '  +0.05%  '
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
